$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "43.107.58"
Set-TextValue $ws "E2" "  +4.50%  "
Set-TextValue $ws "D3" "2.250.64"
Set-TextValue $ws "E3" "  +3.59%  "
Set-TextValue $ws "E4" "  +0.01%  "
Set-TextValue $ws "D5" "245.22"
Set-TextValue $ws "E5" "  +3.60%  "
Set-TextValue $ws "D6" "0.618"
Set-TextValue $ws "E6" "  +1.26%  "
Set-TextValue $ws "D7" "76.41"
Set-TextValue $ws "E7" "  +9.96%  "
Set-TextValue $ws "E8" "  -0.19%  "
Set-TextValue $ws "D9" "0.613"
Set-TextValue $ws "E9" "  +6.59%  "
Set-TextValue $ws "D10" "41.17"
Set-TextValue $ws "E10" "  +4.44%  "
Set-TextValue $ws "D11" "0.0936"
Set-TextValue $ws "E11" "  +2.01%  "
Set-TextValue $ws "D12" "7.01"
Set-TextValue $ws "E12" "  +4.59%  "
Set-TextValue $ws "E13" "  +1.02%  "
Set-TextValue $ws "D14" "2.592.11"
Set-TextValue $ws "E14" "  +3.57%  "
Set-TextValue $ws "D15" "14.66"
Set-TextValue $ws "E15" "  +5.18%  "
Set-TextValue $ws "D16" "2.255.08"
Set-TextValue $ws "E16" "  +3.57%  "
Set-TextValue $ws "E17" "  +1.76%  "
Set-TextValue $ws "D18" "43.022.38"
Set-TextValue $ws "E18" "  +4.78%  "
Set-TextValue $ws "E19" "  +5.67%  "
Set-TextValue $ws "D20" "71.18"
Set-TextValue $ws "E20" "  +0.86%  "
Set-TextValue $ws "E21" "  +2.00%  "
Set-TextValue $ws "D22" "10.11"
Set-TextValue $ws "E22" "  +7.96%  "
Set-TextValue $ws "D23" "231.03"
Set-TextValue $ws "E23" "  +2.49%  "
Set-TextValue $ws "D24" "2.21"
Set-TextValue $ws "E24" "  +16.36%  "
Set-TextValue $ws "E25" "  +0.22%  "
Set-TextValue $ws "D26" "10.93"
Set-TextValue $ws "E26" "  +1.97%  "
Set-TextValue $ws "D27" "3.48"
Set-TextValue $ws "E27" "  -1.22%  "
Set-TextValue $ws "D28" "39.63"
Set-TextValue $ws "E28" "  +29.60%  "
Set-TextValue $ws "E29" "  +3.39%  "
Set-TextValue $ws "E30" "  +1.94%  "
Set-TextValue $ws "D31" "173.71"
Set-TextValue $ws "E31" "  +3.37%  "
Set-TextValue $ws "D32" "20.37"
Set-TextValue $ws "E32" "  +2.42%  "
Set-TextValue $ws "E33" "  +4.93%  "
Set-TextValue $ws "D34" "5.37"
Set-TextValue $ws "E34" "  +5.61%  "
Set-TextValue $ws "E35" "  +1.78%  "
Set-TextValue $ws "D36" "0.111"
Set-TextValue $ws "E36" "  +9.88%  "
Set-TextValue $ws "E37" "  +7.54%  "
Set-TextValue $ws "D38" "0.0336"
Set-TextValue $ws "E38" "  +19.15%  "
Set-TextValue $ws "D39" "13.02"
Set-TextValue $ws "E39" "  +13.14%  "
Set-TextValue $ws "E40" "  +3.97%  "
Set-TextValue $ws "E41" "  +3.04%  "
Set-TextValue $ws "D42" "0.205"
Set-TextValue $ws "E42" "  +8.30%  "
Set-TextValue $ws "D43" "60.17"
Set-TextValue $ws "E43" "  +2.04%  "
Set-TextValue $ws "D44" "106.06"
Set-TextValue $ws "E44" "  +8.85%  "
Set-TextValue $ws "D45" "8.71"
Set-TextValue $ws "E45" "  +5.57%  "
Set-TextValue $ws "D46" "0.1000"
Set-TextValue $ws "E46" "  +3.29%  "
Set-TextValue $ws "D47" "0.464"
Set-TextValue $ws "E47" "  +27.26%  "
Set-TextValue $ws "E48" "  +9.96%  "
Set-TextValue $ws "E49" "  +3.27%  "
Set-TextValue $ws "D50" "1.15"
Set-TextValue $ws "E50" "  +3.07%  "
Set-TextValue $ws "D51" "2.463.80"
Set-TextValue $ws "E51" "  +3.46%  "
